$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 311  # F5: 310 -> 311
$ws.Cells.Item(6, 6).Value = 447  # F6: 445 -> 447
$ws.Cells.Item(8, 6).Value = 2024  # F8: 2018 -> 2024
$ws.Cells.Item(10, 6).Value = 38  # F10: 37 -> 38
$ws.Cells.Item(11, 6).Value = 36  # F11: 35 -> 36
$ws.Cells.Item(12, 6).Value = 1601  # F12: 1600 -> 1601
$ws.Cells.Item(13, 6).Value = 1601  # F13: 1600 -> 1601
$ws.Cells.Item(14, 6).Value = 1331  # F14: 1328 -> 1331
$ws.Cells.Item(16, 6).Value = 1387  # F16: 1386 -> 1387
$ws.Cells.Item(18, 6).Value = 15  # F18: 14 -> 15
$ws.Cells.Item(20, 6).Value = 479  # F20: 474 -> 479
$ws.Cells.Item(22, 6).Value = 147  # F22: 146 -> 147
$ws.Cells.Item(23, 6).Value = 7062  # F23: 7057 -> 7062
$ws.Cells.Item(24, 6).Value = 7062  # F24: 7057 -> 7062
$ws.Cells.Item(25, 6).Value = 7668  # F25: 7648 -> 7668
$ws.Cells.Item(28, 6).Value = 183  # F28: 182 -> 183
$ws.Cells.Item(30, 6).Value = 81  # F30: 80 -> 81
$ws.Cells.Item(32, 6).Value = 250  # F32: 248 -> 250
$ws.Cells.Item(33, 6).Value = 180  # F33: 172 -> 180
$ws.Cells.Item(36, 6).Value = 41  # F36: 40 -> 41
$ws.Cells.Item(38, 6).Value = 1393  # F38: 1391 -> 1393
$ws.Cells.Item(39, 6).Value = 21  # F39: 19 -> 21
$ws.Cells.Item(40, 6).Value = 216  # F40: 215 -> 216
$ws.Cells.Item(41, 6).Value = 283  # F41: 282 -> 283
$ws.Cells.Item(42, 6).Value = 695  # F42: 689 -> 695
$ws.Cells.Item(45, 6).Value = 315  # F45: 313 -> 315
$ws.Cells.Item(47, 6).Value = 188  # F47: 187 -> 188
$ws.Cells.Item(49, 6).Value = 141  # F49: 137 -> 141
$ws.Cells.Item(50, 6).Value = 142  # F50: 140 -> 142

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 7  # F4: 6 -> 7
$ws.Cells.Item(16, 6).Value = 5  # F16: 3 -> 5

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 2583  # F3: 2579 -> 2583
$ws.Cells.Item(4, 6).Value = 264  # F4: 263 -> 264

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 6).Value = 264  # F6: 263 -> 264
$ws.Cells.Item(9, 6).Value = 311  # F9: 310 -> 311
$ws.Cells.Item(11, 6).Value = 447  # F11: 445 -> 447
$ws.Cells.Item(12, 6).Value = 2024  # F12: 2018 -> 2024
$ws.Cells.Item(13, 6).Value = 38  # F13: 37 -> 38
$ws.Cells.Item(14, 6).Value = 36  # F14: 35 -> 36
$ws.Cells.Item(15, 6).Value = 1601  # F15: 1600 -> 1601
$ws.Cells.Item(16, 6).Value = 1601  # F16: 1600 -> 1601
$ws.Cells.Item(17, 6).Value = 1331  # F17: 1328 -> 1331
$ws.Cells.Item(18, 6).Value = 1387  # F18: 1386 -> 1387
$ws.Cells.Item(20, 6).Value = 479  # F20: 474 -> 479
$ws.Cells.Item(21, 6).Value = 7  # F21: 6 -> 7
$ws.Cells.Item(22, 6).Value = 147  # F22: 146 -> 147
$ws.Cells.Item(24, 6).Value = 7062  # F24: 7057 -> 7062
$ws.Cells.Item(25, 6).Value = 7062  # F25: 7057 -> 7062
$ws.Cells.Item(26, 6).Value = 7668  # F26: 7648 -> 7668
$ws.Cells.Item(29, 6).Value = 81  # F29: 80 -> 81
$ws.Cells.Item(30, 6).Value = 250  # F30: 248 -> 250
$ws.Cells.Item(33, 6).Value = 41  # F33: 40 -> 41
$ws.Cells.Item(34, 6).Value = 1393  # F34: 1391 -> 1393
$ws.Cells.Item(35, 6).Value = 21  # F35: 19 -> 21
$ws.Cells.Item(36, 6).Value = 216  # F36: 215 -> 216
$ws.Cells.Item(38, 6).Value = 283  # F38: 282 -> 283
$ws.Cells.Item(41, 6).Value = 695  # F41: 689 -> 695
$ws.Cells.Item(45, 6).Value = 315  # F45: 313 -> 315
$ws.Cells.Item(48, 6).Value = 141  # F48: 137 -> 141
$ws.Cells.Item(49, 6).Value = 5  # F49: 3 -> 5
